$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair data (TPM re-run): columns A..T for rows 2..9
# A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the numeric NATMI metrics.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs", "Wnt1", "Fzd3", "ECs", 1, 0.3333333333333333, 0.131499, 0.394497, 0.3654391092296077, 0.3654391092296077, 3, 1, 0.3643233333333333, 1.09297, 0.09146135066097912, 0.09146135066097912, 0.04790815401, 0.43117338609, 0.033423554514485, 0.033423554514485),
    @("ECs", "Wnt1", "Fzd3", "FAPs", 1, 0.3333333333333333, 0.131499, 0.394497, 0.3654391092296077, 0.3654391092296077, 3, 1, 0.8809360000000001, 2.642808, 0.2211540931751474, 0.2211540931751474, 0.115842203064, 1.042579827576, 0.08081835481240755, 0.08081835481240755),
    @("ECs", "Wnt1", "Fzd3", "MuSCs", 1, 0.3333333333333333, 0.131499, 0.394497, 0.3654391092296077, 0.3654391092296077, 3, 1, 2.721212, 8.163636, 0.6831451685449673, 0.6831451685449673, 0.357836656788, 3.220529911092, 0.2496479618675831, 0.2496479618675831),
    @("ECs", "Wnt1", "Fzd3", "Resolving-Mac", 1, 0.3333333333333333, 0.131499, 0.394497, 0.3654391092296077, 0.3654391092296077, 1, 0.3333333333333333, 0.016887, 0.050661, 0.004239387618906157, 0.004239387618906157, 0.002220623613, 0.019985612517, 0.001549238035132094, 0.001549238035132094),
    @("FAPs", "Wnt1", "Fzd3", "ECs", 1, 0.3333333333333333, 0.2283393333333333, 0.685018, 0.6345608907703922, 0.6345608907703922, 3, 1, 0.3643233333333333, 1.09297, 0.09146135066097912, 0.09146135066097912, 0.08318934705111111, 0.74870412346, 0.05803779614649411, 0.05803779614649411),
    @("FAPs", "Wnt1", "Fzd3", "FAPs", 1, 0.3333333333333333, 0.2283393333333333, 0.685018, 0.6345608907703922, 0.6345608907703922, 3, 1, 0.8809360000000001, 2.642808, 0.2211540931751474, 0.2211540931751474, 0.2011523389493333, 1.810371050544, 0.1403357383627399, 0.1403357383627399),
    @("FAPs", "Wnt1", "Fzd3", "MuSCs", 1, 0.3333333333333333, 0.2283393333333333, 0.685018, 0.6345608907703922, 0.6345608907703922, 3, 1, 2.721212, 8.163636, 0.6831451685449673, 0.6831451685449673, 0.6213597339386667, 5.592237605448, 0.4334972066773842, 0.4334972066773842),
    @("FAPs", "Wnt1", "Fzd3", "Resolving-Mac", 1, 0.3333333333333333, 0.2283393333333333, 0.685018, 0.6345608907703922, 0.6345608907703922, 1, 0.3333333333333333, 0.016887, 0.050661, 0.004239387618906157, 0.004239387618906157, 0.003855966322, 0.034703696898, 0.002690149583774063, 0.002690149583774063)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $columns.Count; $c++) {
        $ws.Range($columns[$c] + $r).Value = $rowData[$c]
    }
}
